$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply every cell value change described by the diff, in row order.
# Cells whose new value happens to look like a plain number (e.g. "208.12",
# "0.490", "1.00" ...) are forced to Text format first -- otherwise the COM
# layer auto-converts the typed value to a real number the way Excel does,
# which would silently drop formatting such as trailing zeros.

$ws.Range("D2").Value = '27.279.33'
$ws.Range("E2").Value = '  -1.40%  '
$ws.Range("D3").Value = '1.576.95'
$ws.Range("E3").Value = '  -0.83%  '
$ws.Range("E4").Value = '  -0.26%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '208.12'
$ws.Range("E5").Value = '  -0.28%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.490'
$ws.Range("E6").Value = '  -1.78%  '
$ws.Range("E7").Value = '  -0.28%  '
$ws.Range("E8").Value = '  +0.14%  '
$ws.Range("E9").Value = '  -1.28%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0591'
$ws.Range("E10").Value = '  +0.21%  '
$ws.Range("E11").Value = '  -0.11%  '
$ws.Range("D12").Value = '1.802.59'
$ws.Range("E12").Value = '  -0.87%  '
$ws.Range("D13").Value = '1.580.68'
$ws.Range("E13").Value = '  -0.41%  '
$ws.Range("E14").Value = '  -1.12%  '
$ws.Range("E15").Value = '  -1.35%  '
$ws.Range("D16").Value = '27.304.65'
$ws.Range("E16").Value = '  -1.41%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '62.54'
$ws.Range("E17").Value = '  -0.95%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '215.34'
$ws.Range("E18").Value = '  -0.90%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '7.35'
$ws.Range("E19").Value = '  +0.27%  '
$ws.Range("D20").Value = '0.0₃0687'
$ws.Range("E20").Value = '  -0.96%  '
$ws.Range("E21").Value = '  -0.22%  '
$ws.Range("E22").Value = '  -0.10%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '9.43'
$ws.Range("E23").Value = '  -3.60%  '
$ws.Range("E24").Value = '  +1.68%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '151.92'
$ws.Range("E25").Value = '  -1.09%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '6.69'
$ws.Range("E26").Value = '  -3.83%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '14.96'
$ws.Range("B28").Value = 'Stellar'
$ws.Range("C28").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.104'
$ws.Range("E28").Value = '  -1.12%  '
$ws.Range("B29").Value = 'BinanceUSD'
$ws.Range("C29").Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.00'
$ws.Range("E29").Value = '  -0.27%  '
$ws.Range("E30").Value = '  -1.53%  '
$ws.Range("E31").Value = '  -1.93%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.18'
$ws.Range("E32").Value = '  -1.16%  '
$ws.Range("D33").Value = '1.411.55'
$ws.Range("E33").Value = '  +2.40%  '
$ws.Range("E34").Value = '  -1.35%  '
$ws.Range("E35").Value = '  +1.70%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.29'
$ws.Range("E36").Value = '  -1.68%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.941'
$ws.Range("E37").Value = '  -2.45%  '
$ws.Range("E38").Value = '  -1.79%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.821'
$ws.Range("E39").Value = '  -0.42%  '
$ws.Range("E41").Value = '  -0.28%  '
$ws.Range("E42").Value = '  +2.04%  '
$ws.Range("E43").Value = '  +3.56%  '
$ws.Range("E44").Value = '  +1.72%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '63.99'
$ws.Range("E45").Value = '  -0.47%  '
$ws.Range("E46").Value = '  +0.36%  '
$ws.Range("D47").Value = '1.714.66'
$ws.Range("E47").Value = '  -0.88%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '86.19'
$ws.Range("E48").Value = '  +0.44%  '
$ws.Range("D49").Value = '0.0₇0984'
$ws.Range("E49").Value = '  -1.56%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0954'
$ws.Range("E50").Value = '  -1.00%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0495'
$ws.Range("E51").Value = '  +0.10%  '
